$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update raw count data (B:D) for rows 2-13 -------------------------
# row -> (NumSNPs, NumMultiallelic, NumInvariants)
$data = @{
    2  = @(165262, 9416, 1259181)
    3  = @(167481, 9231, 1368061)
    4  = @(178061, 11455, 1180563)
    5  = @(130871, 7412, 1002216)
    6  = @(160414, 10443, 981011)
    7  = @(196142, 10318, 1748766)
    8  = @(150019, 9309, 1058566)
    9  = @(236520, 14122, 1654303)
    10 = @(144600, 7985, 1104510)
    11 = @(153849, 8426, 1258706)
    12 = @(149538, 8141, 1262938)
    13 = @(7082, 657, 34194)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
}

# --- New column F: "Without Multiallelic SNPs" -------------------------
# Build the new style (bold font + full box border, general number format)
# by copying the fully-boxed numeric style from E14 and then stripping its
# number format back down to General, so the resulting cell style exactly
# matches: bold font, box border, General format.
$ws.Range("E14").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Value = "Without Multiallelic SNPs"
$ws.Range("F13").Font.Bold = $true
$ws.Range("F13").NumberFormat = "general"

$ws.Range("E14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Formula = "=B14+D14"

# --- Column widths for the new E/F columns ------------------------------
$ws.Columns.Item(5).ColumnWidth = 9.498697916666666
$ws.Columns.Item(6).ColumnWidth = 22.498697916666668

# --- Selection, matching the post-edit cursor position ------------------
[void]$ws.Range("F15").Select()
